$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Carry over the handful of per-cell format quirks that rode along with
# the manual row reshuffle (copy formats first, while the source cells still
# carry their original formatting, then overwrite the cell values). ---

# Row 13's "last row" look (no bottom border, smaller/plain font) moves to
# row 12 (now the new last visually-special row before the final entry).
$ws.Cells.Item(13, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(13, 2).Copy()
$ws.Cells.Item(12, 2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(13, 3).Copy()
$ws.Cells.Item(12, 3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(13, 7).Copy()
$ws.Cells.Item(12, 7).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(13, 8).Copy()
$ws.Cells.Item(12, 8).PasteSpecial($xlPasteFormats)

# Row 13 itself reverts to the regular row formatting (borrow from row 11,
# which keeps the plain style throughout).
$ws.Cells.Item(11, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(10, 7).Copy()
$ws.Cells.Item(13, 7).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(11, 8).Copy()
$ws.Cells.Item(13, 8).PasteSpecial($xlPasteFormats)

# Column C's odd explicit-style/no-style quirk also shifts: row 9 loses it,
# row 11 picks it up.
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(9, 3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 3).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# --- Update the actual server roster / ids / ports (rows 9-13) ---

# Row 9: GameServer_1 / 6 -> MasterServer_1 / 1
$ws.Cells.Item(9, 1).Value = "MasterServer_1"
$ws.Cells.Item(9, 2).Value = "1"
$ws.Cells.Item(9, 3).Value = "MasterServer_1"
$ws.Cells.Item(9, 7).Value = 13001
$ws.Cells.Item(9, 8).Value = "1"

# Row 10: WorldServer_1 / 7 -> WorldServer_1 / 50
$ws.Cells.Item(10, 1).Value = "WorldServer_1"
$ws.Cells.Item(10, 2).Value = "50"
$ws.Cells.Item(10, 3).Value = "WorldServer_1"
$ws.Cells.Item(10, 7).Value = 17001
$ws.Cells.Item(10, 8).Value = "50"

# Row 11: ProxyServer_1 / 5 -> GameServer_1 / 51
$ws.Cells.Item(11, 1).Value = "GameServer_1"
$ws.Cells.Item(11, 2).Value = "51"
$ws.Cells.Item(11, 3).Value = "GameServer_1"
$ws.Cells.Item(11, 7).Value = 16001
$ws.Cells.Item(11, 8).Value = "51"

# Row 12: MasterServer_1 / 3 -> LoginServer_1 / 52
$ws.Cells.Item(12, 1).Value = "LoginServer_1"
$ws.Cells.Item(12, 2).Value = "52"
$ws.Cells.Item(12, 3).Value = "LoginServer_1"
$ws.Cells.Item(12, 7).Value = 14001
$ws.Cells.Item(12, 8).Value = "52"

# Row 13: LoginServer_1 / 4 -> ProxyServer_1 / 53
$ws.Cells.Item(13, 1).Value = "ProxyServer_1"
$ws.Cells.Item(13, 2).Value = "53"
$ws.Cells.Item(13, 3).Value = "ProxyServer_1"
$ws.Cells.Item(13, 7).Value = 15001
$ws.Cells.Item(13, 8).Value = "53"
